# Insert 6 new price-report rows (date 2022-01-27, serial 44588) right before
# what is currently row 306. This shifts the existing rows 306:369 down to
# 312:375, which matches the target workbook exactly (the rows that already
# existed simply move down by 6 positions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("306:311").Insert()

# Columns that are constant for every row in this sheet.
$A = 3
$B = "Femacal de La Calera"
$C = "Coquimbo"
$E = 5
$F = 100112027
$G = "Melón"
$N = "`$/unidad"
$Q = 1
$R = "Hortaliza"

# New data block: date 44588, two varieties (Calameño / Tuna) each with
# Extra / Primera / Segunda qualities, all sourced from Provincia de Quillota.
$newRows = @(
    @{ Row = 306; D = 44588; H = "Calameño"; I = "Extra";   J = 450; K = 1000; L = 1000; M = 1000; O = "Provincia de Quillota"; P = 1000 },
    @{ Row = 307; D = 44588; H = "Calameño"; I = "Primera"; J = 400; K = 700;  L = 700;  M = 700;  O = "Provincia de Quillota"; P = 700  },
    @{ Row = 308; D = 44588; H = "Calameño"; I = "Segunda"; J = 380; K = 500;  L = 500;  M = 500;  O = "Provincia de Quillota"; P = 500  },
    @{ Row = 309; D = 44588; H = "Tuna";     I = "Extra";   J = 280; K = 1000; L = 1000; M = 1000; O = "Provincia de Quillota"; P = 1000 },
    @{ Row = 310; D = 44588; H = "Tuna";     I = "Primera"; J = 300; K = 700;  L = 700;  M = 700;  O = "Provincia de Quillota"; P = 700  },
    @{ Row = 311; D = 44588; H = "Tuna";     I = "Segunda"; J = 260; K = 500;  L = 500;  M = 500;  O = "Provincia de Quillota"; P = 500  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = $A
    $ws.Range("B$row").Value2 = $B
    $ws.Range("C$row").Value2 = $C
    $ws.Range("D$row").Value2 = $r.D
    $ws.Range("E$row").Value2 = $E
    $ws.Range("F$row").Value2 = $F
    $ws.Range("G$row").Value2 = $G
    $ws.Range("H$row").Value2 = $r.H
    $ws.Range("I$row").Value2 = $r.I
    $ws.Range("J$row").Value2 = $r.J
    $ws.Range("K$row").Value2 = $r.K
    $ws.Range("L$row").Value2 = $r.L
    $ws.Range("M$row").Value2 = $r.M
    $ws.Range("N$row").Value2 = $N
    $ws.Range("O$row").Value2 = $r.O
    $ws.Range("P$row").Value2 = $r.P
    $ws.Range("Q$row").Value2 = $Q
    $ws.Range("R$row").Value2 = $R
}
